# Update result values across the scenario-year worksheets with refreshed
# figures from the server (see commit message: "ADD results from server").
$wb = $excel.ActiveWorkbook

# Sheet "2025" (sheet1.xml)
$ws = $wb.Worksheets.Item("2025")
$ws.Range("N2").Value = 5030.117123085906
$ws.Range("O2").Value = 4892.012299288742

# Sheet "2030" (sheet2.xml)
$ws = $wb.Worksheets.Item("2030")
$ws.Range("B2").Value = 3995.471002096464
$ws.Range("I2").Value = 31144.4413139226
$ws.Range("L2").Value = 55661.22279899548
$ws.Range("M2").Value = 15393.99435160624
$ws.Range("N2").Value = 7449.833478787733
$ws.Range("O2").Value = 8471.377716625568

# Sheet "2035" (sheet3.xml)
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 2927.360317916481
$ws.Range("B2").Value = 6228.543249765059
$ws.Range("E2").Value = 50348.17506991202
$ws.Range("I2").Value = 46183.13572784168
$ws.Range("L2").Value = 55661.22279899548
$ws.Range("M2").Value = 18949.69321397918
$ws.Range("N2").Value = 11974.99718757541
$ws.Range("O2").Value = 11153.02281884398

# Sheet "2040" (sheet4.xml)
$ws = $wb.Worksheets.Item("2040")
$ws.Range("A2").Value = 2927.360317916481
$ws.Range("B2").Value = 6228.543249765059
$ws.Range("E2").Value = 50348.17506991202
$ws.Range("I2").Value = 46183.13572784168
$ws.Range("L2").Value = 55661.22279899548
$ws.Range("M2").Value = 18949.69321397918
$ws.Range("N2").Value = 12082.12109559635
$ws.Range("O2").Value = 11153.02281884398

# Sheet "2045" (sheet5.xml)
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 6352.985609279765
$ws.Range("B2").Value = 6228.543249765059
$ws.Range("E2").Value = 50348.17506991202
$ws.Range("I2").Value = 46183.13572784168
$ws.Range("L2").Value = 55661.22279899548
$ws.Range("M2").Value = 18949.69321397918
$ws.Range("N2").Value = 12626.84864704659
$ws.Range("O2").Value = 13488.48880519273

# Sheet "2050" (sheet6.xml)
$ws = $wb.Worksheets.Item("2050")
$ws.Range("A2").Value = 6352.985609279765
$ws.Range("B2").Value = 6228.543249765059
$ws.Range("E2").Value = 50348.17506991202
$ws.Range("I2").Value = 46183.13572784168
$ws.Range("L2").Value = 55661.22279899548
$ws.Range("M2").Value = 18949.69321397918
$ws.Range("N2").Value = 12626.84864704659
$ws.Range("O2").Value = 13488.48880519273
